# "fundamental change of algorithms"
# Quarterly report rolled forward by one quarter:
#   - oldest quarter column ("فصل اول منتهی به 1399/03") is dropped
#   - every quarter's figures shift one column to the left (E<-F, F<-G, ... M<-N)
#   - a brand new quarter ("فصل سوم منتهی به 1401/09") is appended in column N
#   - the two header rows that label the quarters (row 8 and row 24) shift the
#     same way
#   - theme accent1 / accent5 colors are swapped back to the older default
#     Office theme ordering

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Quarter header labels (row 8 and row 24) - shift left, append new quarter
# ---------------------------------------------------------------------------
$quarters = @(
    "فصل دوم منتهی به 1399/06",
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09"
)

$cols = @("E","F","G","H","I","J","K","L","M","N")

foreach ($headerRow in @(8, 24)) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$headerRow"
        $ws.Range($addr).Value2 = $quarters[$i]
    }
}

# ---------------------------------------------------------------------------
# 2. Data rows: shift values one column left, set the new quarter's figure
# ---------------------------------------------------------------------------
$newQuarterValue = @{
    10 = 664256
    13 = 402
    15 = -8671
    16 = 3672
    17 = 30341
    19 = 19440
    20 = 709440
    26 = 579
    27 = 144
}

foreach ($row in $newQuarterValue.Keys) {
    $old = @()
    foreach ($c in $cols) {
        $old += , ($ws.Range("$c$row").Value2)
    }
    # shifted: new E..M = old F..N ; new N = new quarter value
    for ($i = 0; $i -lt ($cols.Length - 1); $i++) {
        $addr = "$($cols[$i])$row"
        $ws.Range($addr).Value2 = $old[$i + 1]
    }
    $ws.Range("N$row").Value2 = $newQuarterValue[$row]
}

# ---------------------------------------------------------------------------
# 3. Theme colors: swap accent1 <-> accent5
# ---------------------------------------------------------------------------
$themeColors = $wb.Theme.ThemeColorScheme
$accent1 = $themeColors.Colors(5)
$accent5 = $themeColors.Colors(9)
$accent1.RGB = 13998939   # 5B9BD5 (BGR-encoded)
$accent5.RGB = 12874308   # 4472C4 (BGR-encoded)
